$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the weight for "Family" (row 2) from 1.0 to 1.2
$ws.Range("B2").Value = 1.2

# Append new communal keywords (Kids, Husband, Wife, Partner) with weight 1.0.
# Column A uses the same style as the header cell (A1); column B uses the
# same style as the other weight cells (B1).
$newKeywords = @("Kids", "Husband", "Wife", "Partner")
$startRow = 51

for ($i = 0; $i -lt $newKeywords.Length; $i++) {
    $row = $startRow + $i

    $ws.Range("A1").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $newKeywords[$i]

    $ws.Range("B1").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
    $ws.Range("B$row").Value = 1.0
}
